# Append the two new daily log rows (2025-11-03 / Excel serial 45964) for the
# "四方坪站" and "高岭站" charging stations onto the bottom of the existing
# table on Sheet1, then move the selection to reflect where the cursor ended
# up after the edit (G134).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 127
$siteFangSquarePing = "四方坪站充电量(kw)"
$siteGaoLing        = "高岭站充电量(kw)"

# Row 1 of 2: 四方坪站, date 45964 (2025-11-03)
$row1 = $lastRow + 1
$values1 = @(
    45964, $siteFangSquarePing,
    720.62400000000002, 846.43000000000018, 446.14000000000004, 488.16600000000005,
    549.76400000000001, 682.923, 598.96400000000006, 180.054,
    93.56, 258.72399999999999, 154.66999999999999, 291.78499999999997,
    739.65300000000002, 1495.8519999999996, 561.79300000000012, 336.68800000000005,
    406.45300000000003, 69.432000000000002, 164.71699999999998, 141.38999999999999,
    140.14000000000001, 24.1, 82.738, 40.43
)

# Row 2 of 2: 高岭站, date 45964 (2025-11-03)
$row2 = $lastRow + 2
$values2 = @(
    45964, $siteGaoLing,
    317.54700000000003, 280.52600000000001, 89.284999999999997, 90.562000000000012,
    113.648, 213.80500000000001, 165.42400000000001, 51.927,
    147.12100000000001, 139.22, 175.637, 220.61899999999997,
    562.68899999999996, 465.57599999999991, 241.18200000000002, 55.242000000000004,
    182.36199999999999, 74.338999999999999, 29.107000000000003, 87.120999999999995,
    7.282, 0, 0, 0
)

for ($col = 1; $col -le $values1.Length; $col++) {
    $ws.Cells.Item($row1, $col).Value2 = $values1[$col - 1]
}

for ($col = 1; $col -le $values2.Length; $col++) {
    $ws.Cells.Item($row2, $col).Value2 = $values2[$col - 1]
}

# Reflect the post-edit active cell/selection.
$ws.Range("G134").Select()
